# Update the two-digit-division worksheet numbers.
#
# The worksheet is a single 20-row x 5-column table; only rows
# 1, 5, 9, 13, 17 (1-based) actually contain equations, the rest are
# blank "work space" rows. We address each cell by (row, column) so the
# edit is unambiguous even though several of the old/new equation
# strings repeat elsewhere in the document.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# row -> array of new cell values (column order 1..5)
$updates = @{
    1  = @("89÷7=", "69÷5=", "59÷9=", "45÷5=", "98÷9=")
    5  = @("74÷7=", "32÷8=", "87÷4=", "49÷7=", "17÷5=")
    9  = @("66÷4=", "26÷4=", "56÷2=", "87÷6=", "17÷5=")
    13 = @("11÷2=", "27÷6=", "90÷5=", "57÷3=", "13÷8=")
    17 = @("82÷9=", "84÷3=", "25÷2=", "97÷6=", "69÷7=")
}

foreach ($rowIndex in $updates.Keys) {
    $values = $updates[$rowIndex]
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
